$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 17:16:38"

# Sheet "zh-cn": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 17:16:33"
$wsZhCn.Range("K2").Value = "2016-09-01 17:16:50"

# Sheet "de-de": Correspond Handoff Datetime (H2) shares the same value as
# Overview!G2 ("Latest HO Xliff Generate Date"), and Correspond Handback
# DateTime (K2) gets its own new timestamp.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 17:16:38"
$wsDeDe.Range("K2").Value = "2016-09-01 17:16:58"
